$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.595.58'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '3.017.47'
$ws.Range("E3").Value = '  +3.33%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '201.43'
$ws.Range("E5").Value = '  +0.95%  '

$ws.Range("D6").Value = '634.27'
$ws.Range("E6").Value = '  +5.80%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  -0.26%  '

$ws.Range("E9").Value = '  +2.76%  '

$ws.Range("D10").Value = '3.021.45'
$ws.Range("E10").Value = '  +3.52%  '

$ws.Range("D11").Value = '0.433'
$ws.Range("E11").Value = '  +1.39%  '

$ws.Range("E12").Value = '  -0.17%  '

$ws.Range("D13").Value = '5.03'
$ws.Range("E13").Value = '  +2.74%  '

$ws.Range("D14").Value = '3.541.19'
$ws.Range("E14").Value = '  +2.55%  '

$ws.Range("D15").Value = '29.22'
$ws.Range("E15").Value = '  +6.51%  '

$ws.Range("D16").Value = '76.422.76'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("D17").Value = '0.0000189'
$ws.Range("E17").Value = '  -1.31%  '

$ws.Range("D18").Value = '2.989.64'
$ws.Range("E18").Value = '  +2.51%  '

$ws.Range("D19").Value = '13.51'
$ws.Range("E19").Value = '  +5.50%  '

$ws.Range("D20").Value = '8.97'
$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("D21").Value = '374.10'
$ws.Range("E21").Value = '  -1.13%  '

$ws.Range("D22").Value = '2.27'
$ws.Range("E22").Value = '  -2.37%  '

$ws.Range("D23").Value = '4.30'
$ws.Range("E23").Value = '  +2.44%  '

$ws.Range("D24").Value = '72.82'
$ws.Range("E24").Value = '  +1.89%  '

$ws.Range("E25").Value = '  +2.21%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").Value = '4.37'
$ws.Range("E27").Value = '  +3.04%  '

$ws.Range("D28").Value = '9.82'
$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("D29").Value = '0.0000107'
$ws.Range("E29").Value = '  -1.82%  '

$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.52%  '

$ws.Range("D31").Value = '8.31'
$ws.Range("E31").Value = '  +7.41%  '

$ws.Range("E32").Value = '  -1.03%  '

$ws.Range("D33").Value = '510.89'
$ws.Range("E33").Value = '  +0.92%  '

$ws.Range("E34").Value = '  +7.63%  '

$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").Value = '20.42'
$ws.Range("E36").Value = '  +1.08%  '

$ws.Range("D37").Value = '163.27'
$ws.Range("E37").Value = '  -0.99%  '

$ws.Range("D38").Value = '19.99'

$ws.Range("D39").Value = '0.383'
$ws.Range("E39").Value = '  +11.60%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '188.04'
$ws.Range("E40").Value = '  +4.20%  '

$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").Value = '0.105'
$ws.Range("E41").Value = '  +10.65%  '

$ws.Range("E42").Value = '  -1.90%  '

$ws.Range("D43").Value = '1.01'
$ws.Range("E43").Value = '  +0.53%  '

$ws.Range("D44").Value = '4.96'
$ws.Range("E44").Value = '  -0.95%  '

$ws.Range("D45").Value = '42.57'
$ws.Range("E45").Value = '  +5.82%  '

$ws.Range("D46").Value = '1.65'
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("D47").Value = '1.24'
$ws.Range("E47").Value = '  +2.47%  '

$ws.Range("D48").Value = '0.710'
$ws.Range("E48").Value = '  +7.34%  '

$ws.Range("D49").Value = '0.595'
$ws.Range("E49").Value = '  +2.97%  '

$ws.Range("D50").Value = '2.32'
$ws.Range("E50").Value = '  -1.08%  '

$ws.Range("D51").Value = '3.84'
$ws.Range("E51").Value = '  +3.27%  '
